# Applies the "Mounted Tire Processing Pipeline" recompute.
# Updates the per-signal-segment distribution (Step1_Data), its running
# cumulative sum (Step2_Sj), and the derived threshold crossing stats
# (Step3_DataPts_0.5 / 0.7 / 0.8 / 0.9) for "signal segment 5" (row 6) and
# "signal segment 10" (row 11) to match the new pipeline output.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$arr = New-Object 'object[,]' 1,33
$arr[0,0] = 0.009721307704564801
$arr[0,1] = 0.2252274311889522
$arr[0,2] = 0.1358976734352286
$arr[0,3] = 0.2201092064371697
$arr[0,4] = 0.01184645294956054
$arr[0,5] = 0.00927755426359452
$arr[0,6] = 0.000507866278389662
$arr[0,7] = 0.1276634423622187
$arr[0,8] = 0.05666894367640326
$arr[0,9] = 0.0535527946875316
$arr[0,10] = 0.003230445561374833
$arr[0,11] = 0.003127123582358345
$arr[0,12] = 0.00659203674404996
$arr[0,13] = 0.001073397691088912
$arr[0,14] = 0.004442970737279658
$arr[0,15] = 0.005961196547886244
$arr[0,16] = 0.0004414292009566197
$arr[0,17] = 0.01592582806184508
$arr[0,18] = 0.006763703888198962
$arr[0,19] = 0.0003411913742581223
$arr[0,20] = 0.02908619340898234
$arr[0,21] = 0.03269095146664273
$arr[0,22] = 0.002496994423325315
$arr[0,23] = 0.01433737111189139
$arr[0,24] = 0.0008645341696735846
$arr[0,25] = 0.006007066036415949
$arr[0,26] = 0.00347717866725051
$arr[0,27] = 0.0005462381450017476
$arr[0,28] = 0.006266396465178174
$arr[0,29] = 0.003115970424207133
$arr[0,30] = 0.002129937957409316
$arr[0,31] = 0.0006091713511115005
$arr[0,32] = 0
$ws.Range("AQ6:BW6").Value = $arr

$arr = New-Object 'object[,]' 1,33
$arr[0,0] = 0.003326798909701525
$arr[0,1] = 0.3027048058516318
$arr[0,2] = 0.1004325187808486
$arr[0,3] = 0.2114559552166512
$arr[0,4] = 0.01600702346025521
$arr[0,5] = 0.00361876103240911
$arr[0,6] = 0.002336507179939995
$arr[0,7] = 0.08061594540617836
$arr[0,8] = 0.02130122639657739
$arr[0,9] = 0.05935150244256566
$arr[0,10] = 0.011153046297353
$arr[0,11] = 0.000027840036626598890339692391
$arr[0,12] = 0.01724680788284936
$arr[0,13] = 0.0129222547017418
$arr[0,14] = 0.03612964869608071
$arr[0,15] = 0.009181589842389377
$arr[0,16] = 0.0007221333426660258
$arr[0,17] = 0.01766641634930232
$arr[0,18] = 0.002977421677598215
$arr[0,19] = 0.002667170101093848
$arr[0,20] = 0.04363807443409427
$arr[0,21] = 0.02312531729232948
$arr[0,22] = 0.0001085838702625002
$arr[0,23] = 0.002312888235635131
$arr[0,24] = 0.0004049667968925057
$arr[0,25] = 0.0008389219996570237
$arr[0,26] = 0.00230482155382527
$arr[0,27] = 0.0002404702740221015
$arr[0,28] = 0.006110486204662043
$arr[0,29] = 0.001199137759229352
$arr[0,30] = 0.003743859909675896
$arr[0,31] = 0.004127098065254554
$arr[0,32] = 0
$ws.Range("AQ11:BW11").Value = $arr

$ws = $wb.Worksheets.Item("Step2_Sj")
$arr = New-Object 'object[,]' 1,33
$arr[0,0] = 0.009721307704564801
$arr[0,1] = 0.2349487388935171
$arr[0,2] = 0.3708464123287457
$arr[0,3] = 0.5909556187659153
$arr[0,4] = 0.6028020717154758
$arr[0,5] = 0.6120796259790704
$arr[0,6] = 0.61258749225746
$arr[0,7] = 0.7402509346196787
$arr[0,8] = 0.796919878296082
$arr[0,9] = 0.8504726729836136
$arr[0,10] = 0.8537031185449885
$arr[0,11] = 0.8568302421273468
$arr[0,12] = 0.8634222788713968
$arr[0,13] = 0.8644956765624857
$arr[0,14] = 0.8689386472997653
$arr[0,15] = 0.8748998438476516
$arr[0,16] = 0.8753412730486082
$arr[0,17] = 0.8912671011104533
$arr[0,18] = 0.8980308049986522
$arr[0,19] = 0.8983719963729103
$arr[0,20] = 0.9274581897818927
$arr[0,21] = 0.9601491412485353
$arr[0,22] = 0.9626461356718606
$arr[0,23] = 0.976983506783752
$arr[0,24] = 0.9778480409534256
$arr[0,25] = 0.9838551069898416
$arr[0,26] = 0.9873322856570921
$arr[0,27] = 0.9878785238020938
$arr[0,28] = 0.994144920267272
$arr[0,29] = 0.9972608906914791
$arr[0,30] = 0.9993908286488884
$arr[0,31] = 0.9999999999999999
$arr[0,32] = 0.9999999999999999
$ws.Range("AQ6:BW6").Value = $arr

$arr = New-Object 'object[,]' 1,33
$arr[0,0] = 0.003326798909701525
$arr[0,1] = 0.3060316047613333
$arr[0,2] = 0.4064641235421819
$arr[0,3] = 0.6179200787588331
$arr[0,4] = 0.6339271022190883
$arr[0,5] = 0.6375458632514974
$arr[0,6] = 0.6398823704314374
$arr[0,7] = 0.7204983158376157
$arr[0,8] = 0.7417995422341931
$arr[0,9] = 0.8011510446767588
$arr[0,10] = 0.8123040909741118
$arr[0,11] = 0.8123319310107384
$arr[0,12] = 0.8295787388935878
$arr[0,13] = 0.8425009935953296
$arr[0,14] = 0.8786306422914103
$arr[0,15] = 0.8878122321337997
$arr[0,16] = 0.8885343654764657
$arr[0,17] = 0.906200781825768
$arr[0,18] = 0.9091782035033662
$arr[0,19] = 0.9118453736044601
$arr[0,20] = 0.9554834480385543
$arr[0,21] = 0.9786087653308838
$arr[0,22] = 0.9787173492011463
$arr[0,23] = 0.9810302374367814
$arr[0,24] = 0.9814352042336739
$arr[0,25] = 0.982274126233331
$arr[0,26] = 0.9845789477871563
$arr[0,27] = 0.9848194180611783
$arr[0,28] = 0.9909299042658404
$arr[0,29] = 0.9921290420250697
$arr[0,30] = 0.9958729019347455
$arr[0,31] = 1
$arr[0,32] = 1
$ws.Range("AQ11:BW11").Value = $arr

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F6").Value = 0.5909556187659153
$ws.Range("F11").Value = 0.6179200787588331

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F6").Value = 0.7402509346196787
$ws.Range("D11").Value = 49
$ws.Range("F11").Value = 0.7204983158376157
$ws.Range("G11").Value = 8

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F6").Value = 0.8504726729836136
$ws.Range("D11").Value = 51
$ws.Range("F11").Value = 0.8011510446767588
$ws.Range("G11").Value = 10

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F6").Value = 0.9274581897818927
$ws.Range("D11").Value = 59
$ws.Range("F11").Value = 0.906200781825768
$ws.Range("G11").Value = 18

Write-Output "Mounted Tire Processing Pipeline update applied."
